$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking strings (e.g. "1.030", "6.005")
# keep their exact textual representation instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.617.33"
$ws.Range("E2").Value = "  +1.46%  "

$ws.Range("D3").Value = "1.881.68"
$ws.Range("E3").Value = "  +0.04%  "

$ws.Range("D4").Value = "1.030"
$ws.Range("E4").Value = "  +2.24%  "

$ws.Range("D5").Value = "319.19"
$ws.Range("E5").Value = "  +1.46%  "

$ws.Range("D6").Value = "1.027"
$ws.Range("E6").Value = "  +1.98%  "

$ws.Range("D7").Value = "0.5163"
$ws.Range("E7").Value = "  +0.40%  "

$ws.Range("D8").Value = "0.3957"
$ws.Range("E8").Value = "  +0.91%  "

$ws.Range("E9").Value = "  -0.28%  "

$ws.Range("D10").Value = "1.120"
$ws.Range("E10").Value = "  -0.16%  "

$ws.Range("D11").Value = "42.27"
$ws.Range("E11").Value = "  +1.61%  "

$ws.Range("D12").Value = "6.288"
$ws.Range("E12").Value = "  +0.85%  "

$ws.Range("D13").Value = "20.61"
$ws.Range("E13").Value = "  -0.33%  "

$ws.Range("D14").Value = "1.848.17"
$ws.Range("E14").Value = "  -2.20%  "

$ws.Range("D15").Value = "1.031"
$ws.Range("E15").Value = "  +2.34%  "

$ws.Range("D16").Value = "7.264"
$ws.Range("E16").Value = "  -0.05%  "

$ws.Range("D17").Value = "0.00001113"
$ws.Range("E17").Value = "  +0.99%  "

$ws.Range("D18").Value = "91.62"
$ws.Range("E18").Value = "  +0.45%  "

$ws.Range("D19").Value = "0.06799"
$ws.Range("E19").Value = "  +1.77%  "

$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "1.027"
$ws.Range("E20").Value = "  +2.05%  "

$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "17.76"
$ws.Range("E21").Value = "  -0.30%  "

$ws.Range("D22").Value = "6.005"
$ws.Range("E22").Value = "  -0.80%  "

$ws.Range("D23").Value = "28.659.60"
$ws.Range("E23").Value = "  +1.49%  "

$ws.Range("D24").Value = "11.20"
$ws.Range("E24").Value = "  +0.39%  "

$ws.Range("D25").Value = "2.281"
$ws.Range("E25").Value = "  +0.73%  "

$ws.Range("D26").Value = "162.56"
$ws.Range("E26").Value = "  +1.55%  "

$ws.Range("D27").Value = "2.059.38"
$ws.Range("E27").Value = "  -1.66%  "

$ws.Range("D28").Value = "20.87"

$ws.Range("D29").Value = "2.386"
$ws.Range("E29").Value = "  -4.00%  "

$ws.Range("D30").Value = "127.95"
$ws.Range("E30").Value = "  +2.10%  "

$ws.Range("D31").Value = "0.1056"
$ws.Range("E31").Value = "  -0.60%  "

$ws.Range("D32").Value = "1.041"
$ws.Range("E32").Value = "  +0.10%  "

$ws.Range("D33").Value = "5.854"
$ws.Range("E33").Value = "  +0.07%  "

$ws.Range("D34").Value = "3.665"
$ws.Range("E34").Value = "  +1.37%  "

$ws.Range("D35").Value = "0.02441"
$ws.Range("E35").Value = "  -0.45%  "

$ws.Range("E36").Value = "  -0.51%  "

$ws.Range("D37").Value = "9.216"
$ws.Range("E37").Value = "  -4.68%  "

$ws.Range("D38").Value = "0.2187"
$ws.Range("E38").Value = "  -0.13%  "

$ws.Range("D39").Value = "1.254"
$ws.Range("E39").Value = "  +2.42%  "

$ws.Range("D40").Value = "0.6474"
$ws.Range("E40").Value = "  -0.42%  "

$ws.Range("D41").Value = "1.192"
$ws.Range("E41").Value = "  -0.82%  "

$ws.Range("D42").Value = "5.016"
$ws.Range("E42").Value = "  +0.32%  "

$ws.Range("E43").Value = "  -0.98%  "

$ws.Range("D44").Value = "0.6063"
$ws.Range("E44").Value = "  -1.35%  "

$ws.Range("D45").Value = "13.12"
$ws.Range("E45").Value = "  +0.39%  "

$ws.Range("D46").Value = "3.725"
$ws.Range("E46").Value = "  +1.24%  "

$ws.Range("D47").Value = "1.251"
$ws.Range("E47").Value = "  -2.77%  "

$ws.Range("D48").Value = "2.005"
$ws.Range("E48").Value = "  -0.88%  "

$ws.Range("D49").Value = "1.216"
$ws.Range("E49").Value = "  -1.19%  "

$ws.Range("D50").Value = "122.29"
$ws.Range("E50").Value = "  +1.18%  "

$ws.Range("D51").Value = "0.06881"
$ws.Range("E51").Value = "  -0.52%  "
